$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filter save games) - update B, C, D, E, G for rows 2-21
$ws.Cells.Item(2, 2).Value = 3.182878228561681
$ws.Cells.Item(2, 3).Value = 1.65323645889881
$ws.Cells.Item(2, 4).Value = 0.7127328510149897
$ws.Cells.Item(2, 5).Value = 0.4998867070740569
$ws.Cells.Item(2, 7).Value = 6.048734245549538
$ws.Cells.Item(3, 2).Value = 3.182878228561681
$ws.Cells.Item(3, 3).Value = 1.65323645889881
$ws.Cells.Item(3, 4).Value = 0.7127328510149897
$ws.Cells.Item(3, 5).Value = 0.4998867070740569
$ws.Cells.Item(3, 7).Value = 6.048734245549538
$ws.Cells.Item(4, 2).Value = 3.182878228561681
$ws.Cells.Item(4, 3).Value = 1.65323645889881
$ws.Cells.Item(4, 4).Value = 0.7127328510149897
$ws.Cells.Item(4, 5).Value = 0.4998867070740569
$ws.Cells.Item(4, 7).Value = 6.048734245549538
$ws.Cells.Item(5, 2).Value = 1.505614041169197
$ws.Cells.Item(5, 3).Value = 1.65323645889881
$ws.Cells.Item(5, 4).Value = 0.1529057820181812
$ws.Cells.Item(5, 5).Value = 0.4998867070740569
$ws.Cells.Item(5, 7).Value = 3.811642989160245
$ws.Cells.Item(6, 2).Value = 1.505614041169197
$ws.Cells.Item(6, 3).Value = 1.65323645889881
$ws.Cells.Item(6, 4).Value = 0.7127328510149897
$ws.Cells.Item(6, 5).Value = 0.4998867070740569
$ws.Cells.Item(6, 7).Value = 4.371470058157054
$ws.Cells.Item(7, 2).Value = 1.505614041169197
$ws.Cells.Item(7, 3).Value = 1.65323645889881
$ws.Cells.Item(7, 4).Value = 0.1529057820181812
$ws.Cells.Item(7, 5).Value = 0.4998867070740569
$ws.Cells.Item(7, 7).Value = 3.811642989160245
$ws.Cells.Item(8, 2).Value = 0.7287194209349384
$ws.Cells.Item(8, 3).Value = 1.65323645889881
$ws.Cells.Item(8, 4).Value = 0.1529057820181812
$ws.Cells.Item(8, 5).Value = 0.4998867070740569
$ws.Cells.Item(8, 7).Value = 3.034748368925986
$ws.Cells.Item(9, 2).Value = 3.182878228561681
$ws.Cells.Item(9, 3).Value = 1.65323645889881
$ws.Cells.Item(9, 4).Value = 0.7127328510149897
$ws.Cells.Item(9, 5).Value = 0.4998867070740569
$ws.Cells.Item(9, 7).Value = 6.048734245549538
$ws.Cells.Item(10, 2).Value = 3.182878228561681
$ws.Cells.Item(10, 3).Value = 1.65323645889881
$ws.Cells.Item(10, 4).Value = 0.1529057820181812
$ws.Cells.Item(10, 5).Value = 0.4998867070740569
$ws.Cells.Item(10, 7).Value = 5.488907176552729
$ws.Cells.Item(11, 2).Value = 3.182878228561681
$ws.Cells.Item(11, 3).Value = 1.65323645889881
$ws.Cells.Item(11, 4).Value = 0.1529057820181812
$ws.Cells.Item(11, 5).Value = 0.4998867070740569
$ws.Cells.Item(11, 7).Value = 5.488907176552729
$ws.Cells.Item(12, 2).Value = 3.182878228561681
$ws.Cells.Item(12, 3).Value = 0.3375848360084654
$ws.Cells.Item(12, 4).Value = 0.1529057820181812
$ws.Cells.Item(12, 5).Value = 0.4998867070740569
$ws.Cells.Item(12, 7).Value = 4.173255553662385
$ws.Cells.Item(13, 2).Value = 0.7287194209349384
$ws.Cells.Item(13, 3).Value = 0.3375848360084654
$ws.Cells.Item(13, 4).Value = 0.7127328510149897
$ws.Cells.Item(13, 5).Value = 6.48142807727062
$ws.Cells.Item(13, 7).Value = 8.260465185229014
$ws.Cells.Item(14, 2).Value = 0.7287194209349384
$ws.Cells.Item(14, 3).Value = 1.65323645889881
$ws.Cells.Item(14, 4).Value = 0.1529057820181812
$ws.Cells.Item(14, 5).Value = 0.4998867070740569
$ws.Cells.Item(14, 7).Value = 3.034748368925986
$ws.Cells.Item(15, 2).Value = 1.505614041169197
$ws.Cells.Item(15, 3).Value = 1.65323645889881
$ws.Cells.Item(15, 4).Value = 3.082599426703578
$ws.Cells.Item(15, 5).Value = 6.48142807727062
$ws.Cells.Item(15, 7).Value = 12.7228780040422
$ws.Cells.Item(16, 2).Value = 0.7287194209349384
$ws.Cells.Item(16, 3).Value = 0.3375848360084654
$ws.Cells.Item(16, 4).Value = 0.7127328510149897
$ws.Cells.Item(16, 5).Value = 0.4998867070740569
$ws.Cells.Item(16, 7).Value = 2.27892381503245
$ws.Cells.Item(17, 2).Value = 1.505614041169197
$ws.Cells.Item(17, 3).Value = 1.65323645889881
$ws.Cells.Item(17, 4).Value = 0.1529057820181812
$ws.Cells.Item(17, 5).Value = 0.4998867070740569
$ws.Cells.Item(17, 7).Value = 3.811642989160245
$ws.Cells.Item(18, 2).Value = 0.1554434735375247
$ws.Cells.Item(18, 3).Value = 86.29678392075563
$ws.Cells.Item(18, 4).Value = 0.1529057820181812
$ws.Cells.Item(18, 5).Value = 6.48142807727062
$ws.Cells.Item(18, 7).Value = 93.08656125358196
$ws.Cells.Item(19, 2).Value = 3.182878228561681
$ws.Cells.Item(19, 3).Value = 1.65323645889881
$ws.Cells.Item(19, 4).Value = 3.082599426703578
$ws.Cells.Item(19, 5).Value = 0.4998867070740569
$ws.Cells.Item(19, 7).Value = 8.418600821238126
$ws.Cells.Item(20, 2).Value = 3.182878228561681
$ws.Cells.Item(20, 3).Value = 1.65323645889881
$ws.Cells.Item(20, 4).Value = 3.082599426703578
$ws.Cells.Item(20, 5).Value = 0.4998867070740569
$ws.Cells.Item(20, 7).Value = 8.418600821238126
$ws.Cells.Item(21, 2).Value = 0.06328177979961902
$ws.Cells.Item(21, 3).Value = 0.3375848360084654
$ws.Cells.Item(21, 4).Value = 0.1529057820181812
$ws.Cells.Item(21, 5).Value = 0.4998867070740569
$ws.Cells.Item(21, 7).Value = 1.053659104900323
